# edit.ps1 - apply the TestPlan.docx revisions described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the text spanned by [rangeStart, rangeEnd) with a
# hand-built OOXML fragment (runs / proofErr marks / etc.) while leaving
# the rest of the paragraph (pPr, paraId, rsids, ...) untouched. Must
# re-fetch a *fresh* Range via $d.Range(start, end) - reusing the Find
# range object itself causes InsertXML to append instead of replace.
# ---------------------------------------------------------------------
function Set-RangeXml($rangeStart, $rangeEnd, $innerXml) {
    $pkgHead = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $pkgTail = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pkg = $pkgHead + $innerXml + $pkgTail
    $target = $d.Range($rangeStart, $rangeEnd)
    $target.InsertXML($pkg)
}

# 1. Scope paragraph: drop the comma before "and employee profile viewing."
$d.Content.Find.Execute(
    "leave requests, and employee profile viewing.", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "leave requests and employee profile viewing.", 2) | Out-Null

# 2. Resources bullet: "Test environment (hardware, software)" -> "(software)"
$d.Content.Find.Execute(
    "Test environment (hardware, software)", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Test environment (software)", 2) | Out-Null

# 3. Resources bullet: "Communication tools for issue reporting"
#    -> "Communication tools for issuing a leave." split across two runs.
$rng = $d.Content
$rng.Find.Execute(
    "Communication tools for issue reporting", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml3 = '<w:r><w:t>Communication tools for issu</w:t></w:r>' +
        '<w:r><w:t>ing a leave.</w:t></w:r>'
Set-RangeXml $rng.Start $rng.End $xml3

# 4. "UC-02: Searching Employee for their Department" -> split the run and
#    wrap "Department" in gramStart/gramEnd proofErr marks.
$rng = $d.Content
$rng.Find.Execute(
    "UC-02: Searching Employee for their Department", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml4 = '<w:r><w:t xml:space="preserve">UC-02: Searching Employee for their </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>Department</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>'
Set-RangeXml $rng.Start $rng.End $xml4

# 5. UC-10 bullet: wrap "first time" in gramStart/gramEnd proofErr marks.
$rng = $d.Content
$rng.Find.Execute(
    "Verify that employees can log in for the first time using auto-generated password sent through email.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml5 = '<w:r><w:t xml:space="preserve">Verify that employees can log in for the </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>first time</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> using auto-generated password sent through email.</w:t></w:r>'
Set-RangeXml $rng.Start $rng.End $xml5

# 6. UC-10 bullet: wrap "incorrect" in gramStart/gramEnd proofErr marks.
$rng = $d.Content
$rng.Find.Execute(
    "Verify handling of scenarios where the input information is incorrect or the new password is not valid.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml6 = '<w:r><w:t xml:space="preserve">Verify handling of scenarios where the input information is </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>incorrect</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> or the new password is not valid.</w:t></w:r>'
Set-RangeXml $rng.Start $rng.End $xml6
